# Auto-generated Excel COM-interop script
# Updates market-price / profit columns (H:N) across multiple sheets
# per the scheduled-runner refresh described in the commit diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 436.85715  # H28
$ws.Cells.Item(28, 9).Value = 763.3333  # I28
$ws.Cells.Item(28, 10).Value = 192  # J28
$ws.Cells.Item(28, 11).Value = 763.3333  # K28
$ws.Cells.Item(28, 12).Value = 192  # L28
$ws.Cells.Item(28, 13).Value = -278.3333  # M28
$ws.Cells.Item(28, 14).Value = -1162  # N28
$ws.Cells.Item(74, 8).Value = 1495.7  # H74
$ws.Cells.Item(74, 9).Value = 1494.4155  # I74
$ws.Cells.Item(74, 11).Value = 1494.4155  # K74
$ws.Cells.Item(74, 13).Value = -558.4155000000001  # M74
$ws.Cells.Item(77, 8).Value = 1495.7  # H77
$ws.Cells.Item(77, 9).Value = 1494.4155  # I77
$ws.Cells.Item(77, 11).Value = 7472.0775  # K77
$ws.Cells.Item(77, 13).Value = -2792.0775  # M77
$ws.Cells.Item(99, 8).Value = 364  # H99
$ws.Cells.Item(99, 10).Value = 0  # J99
$ws.Cells.Item(99, 12).Value = 0  # L99
$ws.Cells.Item(99, 14).ClearContents()  # N99
$ws.Cells.Item(129, 8).Value = 1038.3877  # H129
$ws.Cells.Item(129, 10).Value = 1039.1875  # J129
$ws.Cells.Item(129, 12).Value = 3117.5625  # L129
$ws.Cells.Item(129, 14).Value = -13117.5625  # N129
$ws.Cells.Item(132, 8).Value = 49057.227  # H132
$ws.Cells.Item(132, 9).Value = 49057.227  # I132
$ws.Cells.Item(132, 11).Value = 147171.681  # K132
$ws.Cells.Item(132, 13).Value = -144641.681  # M132
$ws.Cells.Item(138, 8).Value = 2470.8518  # H138
$ws.Cells.Item(138, 9).Value = 1940.2142  # I138
$ws.Cells.Item(138, 10).Value = 3042.3076  # J138
$ws.Cells.Item(138, 11).Value = 5820.642599999999  # K138
$ws.Cells.Item(138, 12).Value = 9126.9228  # L138
$ws.Cells.Item(138, 13).Value = -680.6425999999992  # M138
$ws.Cells.Item(138, 14).Value = -19406.9228  # N138
$ws.Cells.Item(139, 8).Value = 52585  # H139
$ws.Cells.Item(139, 10).Value = 52585  # J139
$ws.Cells.Item(139, 12).Value = 52585  # L139
$ws.Cells.Item(139, 14).Value = -62865  # N139

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(122, 8).Value = 2100.3157  # H122
$ws.Cells.Item(122, 9).Value = 1945.2759  # I122
$ws.Cells.Item(122, 11).Value = 5835.8277  # K122
$ws.Cells.Item(122, 13).Value = -3385.8277  # M122
$ws.Cells.Item(135, 8).Value = 28704.834  # H135
$ws.Cells.Item(135, 10).Value = 28704.834  # J135
$ws.Cells.Item(135, 12).Value = 28704.834  # L135
$ws.Cells.Item(135, 14).Value = -38844.834  # N135

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 1133.3334  # H20
$ws.Cells.Item(20, 9).Value = 1150  # I20
$ws.Cells.Item(20, 11).Value = 1150  # K20
$ws.Cells.Item(20, 13).Value = -903  # M20
$ws.Cells.Item(107, 8).Value = 1552.75  # H107
$ws.Cells.Item(107, 9).Value = 1605.5  # I107
$ws.Cells.Item(107, 11).Value = 1605.5  # K107
$ws.Cells.Item(107, 13).Value = 314.5  # M107
$ws.Cells.Item(134, 8).Value = 3826.0417  # H134
$ws.Cells.Item(134, 9).Value = 4349.263  # I134
$ws.Cells.Item(134, 11).Value = 13047.789  # K134
$ws.Cells.Item(134, 13).Value = -10512.789  # M134

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 15570.24  # H31
$ws.Cells.Item(31, 9).Value = 51518.832  # I31
$ws.Cells.Item(31, 10).Value = 4218.0527  # J31
$ws.Cells.Item(31, 11).Value = 51518.832  # K31
$ws.Cells.Item(31, 12).Value = 4218.0527  # L31
$ws.Cells.Item(31, 13).Value = -51223.832  # M31
$ws.Cells.Item(31, 14).Value = -4808.0527  # N31
$ws.Cells.Item(34, 8).Value = 15570.24  # H34
$ws.Cells.Item(34, 9).Value = 51518.832  # I34
$ws.Cells.Item(34, 10).Value = 4218.0527  # J34
$ws.Cells.Item(34, 11).Value = 51518.832  # K34
$ws.Cells.Item(34, 12).Value = 4218.0527  # L34
$ws.Cells.Item(34, 13).Value = -51316.832  # M34
$ws.Cells.Item(34, 14).Value = -4622.0527  # N34
$ws.Cells.Item(58, 8).Value = 13187.634  # H58
$ws.Cells.Item(58, 9).Value = 1145.6552  # I58
$ws.Cells.Item(58, 10).Value = 42289.082  # J58
$ws.Cells.Item(58, 11).Value = 1145.6552  # K58
$ws.Cells.Item(58, 12).Value = 42289.082  # L58
$ws.Cells.Item(58, 13).Value = -942.6551999999999  # M58
$ws.Cells.Item(58, 14).Value = -42695.082  # N58
$ws.Cells.Item(94, 8).Value = 3335.7058  # H94
$ws.Cells.Item(94, 9).Value = 1744.1111  # I94
$ws.Cells.Item(94, 10).Value = 5126.25  # J94
$ws.Cells.Item(94, 11).Value = 1744.1111  # K94
$ws.Cells.Item(94, 12).Value = 5126.25  # L94
$ws.Cells.Item(94, 13).Value = -1293.1111  # M94
$ws.Cells.Item(94, 14).Value = -6028.25  # N94
$ws.Cells.Item(99, 8).Value = 5605.263  # H99
$ws.Cells.Item(99, 9).Value = 4150  # I99
$ws.Cells.Item(99, 10).Value = 7222.222  # J99
$ws.Cells.Item(99, 11).Value = 4150  # K99
$ws.Cells.Item(99, 12).Value = 7222.222  # L99
$ws.Cells.Item(99, 13).Value = -2652  # M99
$ws.Cells.Item(99, 14).Value = -10218.222  # N99
$ws.Cells.Item(107, 8).Value = 416.45456  # H107
$ws.Cells.Item(107, 10).Value = 233.33333  # J107
$ws.Cells.Item(107, 12).Value = 233.33333  # L107
$ws.Cells.Item(107, 14).Value = -4073.33333  # N107
$ws.Cells.Item(126, 8).Value = 5605.263  # H126
$ws.Cells.Item(126, 9).Value = 4150  # I126
$ws.Cells.Item(126, 10).Value = 7222.222  # J126
$ws.Cells.Item(126, 11).Value = 12450  # K126
$ws.Cells.Item(126, 12).Value = 21666.666  # L126
$ws.Cells.Item(126, 13).Value = -9980  # M126
$ws.Cells.Item(126, 14).Value = -26606.666  # N126
$ws.Cells.Item(132, 8).Value = 23001.84  # H132
$ws.Cells.Item(132, 9).Value = 51371.3  # I132
$ws.Cells.Item(132, 11).Value = 154113.9  # K132
$ws.Cells.Item(132, 13).Value = -151583.9  # M132
$ws.Cells.Item(136, 8).Value = 13187.634  # H136
$ws.Cells.Item(136, 9).Value = 1145.6552  # I136
$ws.Cells.Item(136, 10).Value = 42289.082  # J136
$ws.Cells.Item(136, 11).Value = 3436.9656  # K136
$ws.Cells.Item(136, 12).Value = 126867.246  # L136
$ws.Cells.Item(136, 13).Value = -886.9655999999995  # M136
$ws.Cells.Item(136, 14).Value = -131967.246  # N136

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(3, 8).Value = 5791.6  # H3
$ws.Cells.Item(3, 10).Value = 7986  # J3
$ws.Cells.Item(3, 12).Value = 23958  # L3
$ws.Cells.Item(3, 14).Value = -24182  # N3
$ws.Cells.Item(23, 8).Value = 679.8125  # H23
$ws.Cells.Item(23, 9).Value = 133.66667  # I23
$ws.Cells.Item(23, 10).Value = 805.8461  # J23
$ws.Cells.Item(23, 11).Value = 401.00001  # K23
$ws.Cells.Item(23, 12).Value = 2417.5383  # L23
$ws.Cells.Item(23, 13).Value = -166.00001  # M23
$ws.Cells.Item(23, 14).Value = -2887.5383  # N23
$ws.Cells.Item(39, 8).Value = 2254.2  # H39
$ws.Cells.Item(39, 10).Value = 2254.2  # J39
$ws.Cells.Item(39, 12).Value = 6762.599999999999  # L39
$ws.Cells.Item(39, 14).Value = -7350.599999999999  # N39
$ws.Cells.Item(52, 8).Value = 388.41666  # H52
$ws.Cells.Item(52, 10).Value = 388.41666  # J52
$ws.Cells.Item(52, 12).Value = 1165.24998  # L52
$ws.Cells.Item(52, 14).Value = -1697.24998  # N52
$ws.Cells.Item(63, 8).Value = 4015.5  # H63
$ws.Cells.Item(63, 9).Value = 2680  # I63
$ws.Cells.Item(63, 10).Value = 5351  # J63
$ws.Cells.Item(63, 11).Value = 8040  # K63
$ws.Cells.Item(63, 12).Value = 16053  # L63
$ws.Cells.Item(63, 13).Value = -7291  # M63
$ws.Cells.Item(63, 14).Value = -17551  # N63
$ws.Cells.Item(64, 8).Value = 1468.5  # H64
$ws.Cells.Item(64, 9).Value = 1076.1111  # I64
$ws.Cells.Item(64, 10).Value = 5000  # J64
$ws.Cells.Item(64, 11).Value = 3228.3333  # K64
$ws.Cells.Item(64, 12).Value = 15000  # L64
$ws.Cells.Item(64, 13).Value = -2958.3333  # M64
$ws.Cells.Item(64, 14).Value = -15540  # N64
$ws.Cells.Item(66, 8).Value = 4015.5  # H66
$ws.Cells.Item(66, 9).Value = 2680  # I66
$ws.Cells.Item(66, 10).Value = 5351  # J66
$ws.Cells.Item(66, 11).Value = 24120  # K66
$ws.Cells.Item(66, 12).Value = 48159  # L66
$ws.Cells.Item(66, 13).Value = -20376  # M66
$ws.Cells.Item(66, 14).Value = -55647  # N66
$ws.Cells.Item(67, 8).Value = 1468.5  # H67
$ws.Cells.Item(67, 9).Value = 1076.1111  # I67
$ws.Cells.Item(67, 10).Value = 5000  # J67
$ws.Cells.Item(67, 11).Value = 3228.3333  # K67
$ws.Cells.Item(67, 12).Value = 15000  # L67
$ws.Cells.Item(67, 13).Value = -2292.3333  # M67
$ws.Cells.Item(67, 14).Value = -16872  # N67
$ws.Cells.Item(68, 8).Value = 1386.6207  # H68
$ws.Cells.Item(68, 10).Value = 1646.5  # J68
$ws.Cells.Item(68, 12).Value = 4939.5  # L68
$ws.Cells.Item(68, 14).Value = -6561.5  # N68
$ws.Cells.Item(71, 8).Value = 1386.6207  # H71
$ws.Cells.Item(71, 10).Value = 1646.5  # J71
$ws.Cells.Item(71, 12).Value = 14818.5  # L71
$ws.Cells.Item(71, 14).Value = -22930.5  # N71
$ws.Cells.Item(93, 8).Value = 3550  # H93
$ws.Cells.Item(93, 9).Value = 3700  # I93
$ws.Cells.Item(93, 11).Value = 11100  # K93
$ws.Cells.Item(93, 13).Value = -9228  # M93
$ws.Cells.Item(103, 8).Value = 1494.6666  # H103
$ws.Cells.Item(103, 9).Value = 630.44446  # I103
$ws.Cells.Item(103, 10).Value = 2791  # J103
$ws.Cells.Item(103, 11).Value = 1891.33338  # K103
$ws.Cells.Item(103, 12).Value = 8373  # L103
$ws.Cells.Item(103, 13).Value = -1012.33338  # M103
$ws.Cells.Item(103, 14).Value = -10131  # N103
$ws.Cells.Item(131, 8).Value = 774.36  # H131
$ws.Cells.Item(131, 10).Value = 803.4316  # J131
$ws.Cells.Item(131, 12).Value = 2410.2948  # L131
$ws.Cells.Item(131, 14).Value = -12490.2948  # N131
$ws.Cells.Item(137, 8).Value = 7829.227  # H137
$ws.Cells.Item(137, 9).Value = 34666  # I137
$ws.Cells.Item(137, 10).Value = 3591.842  # J137
$ws.Cells.Item(137, 11).Value = 103998  # K137
$ws.Cells.Item(137, 12).Value = 10775.526  # L137
$ws.Cells.Item(137, 13).Value = -98898  # M137
$ws.Cells.Item(137, 14).Value = -20975.526  # N137

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 4449.875  # H70
$ws.Cells.Item(70, 9).Value = 4266.6665  # I70
$ws.Cells.Item(70, 10).Value = 4559.8  # J70
$ws.Cells.Item(70, 11).Value = 4266.6665  # K70
$ws.Cells.Item(70, 12).Value = 4559.8  # L70
$ws.Cells.Item(70, 13).Value = -3996.6665  # M70
$ws.Cells.Item(70, 14).Value = -5099.8  # N70
$ws.Cells.Item(73, 8).Value = 4449.875  # H73
$ws.Cells.Item(73, 9).Value = 4266.6665  # I73
$ws.Cells.Item(73, 10).Value = 4559.8  # J73
$ws.Cells.Item(73, 11).Value = 4266.6665  # K73
$ws.Cells.Item(73, 12).Value = 4559.8  # L73
$ws.Cells.Item(73, 13).Value = -3330.6665  # M73
$ws.Cells.Item(73, 14).Value = -6431.8  # N73
$ws.Cells.Item(113, 8).Value = 3326.5186  # H113
$ws.Cells.Item(113, 9).Value = 2825.85  # I113
$ws.Cells.Item(113, 10).Value = 4757  # J113
$ws.Cells.Item(113, 11).Value = 2825.85  # K113
$ws.Cells.Item(113, 12).Value = 4757  # L113
$ws.Cells.Item(113, 13).Value = -655.8499999999999  # M113
$ws.Cells.Item(113, 14).Value = -9097  # N113

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value = 34121.938  # H136
$ws.Cells.Item(136, 9).Value = 64274.625  # I136
$ws.Cells.Item(136, 10).Value = 3969.25  # J136
$ws.Cells.Item(136, 11).Value = 192823.875  # K136
$ws.Cells.Item(136, 12).Value = 11907.75  # L136
$ws.Cells.Item(136, 13).Value = -190273.875  # M136
$ws.Cells.Item(136, 14).Value = -17007.75  # N136

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(46, 8).Value = 28000  # H46
$ws.Cells.Item(46, 10).Value = 28000  # J46
$ws.Cells.Item(46, 12).Value = 28000  # L46
$ws.Cells.Item(46, 14).Value = -28462  # N46
$ws.Cells.Item(134, 8).Value = 28000  # H134
$ws.Cells.Item(134, 10).Value = 28000  # J134
$ws.Cells.Item(134, 12).Value = 84000  # L134
$ws.Cells.Item(134, 14).Value = -89070  # N134
$ws.Cells.Item(136, 8).Value = 1280.6666  # H136
$ws.Cells.Item(136, 9).Value = 804.875  # I136
$ws.Cells.Item(136, 11).Value = 2414.625  # K136
$ws.Cells.Item(136, 13).Value = 135.375  # M136
